$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- TextBox 58 (shape index 22): add "(ChatCompletionAgent)" paragraph and resize ---
$shp = $s.Shapes.Item(22)
$tr = $shp.TextFrame.TextRange
$oldLen = $tr.Length
$newText = "(ChatCompletionAgent)"
$tr.InsertAfter("`r" + $newText) | Out-Null
$startOfNew = $oldLen + 2
$newPart = $tr.Characters($startOfNew, $newText.Length)
$newPart.Font.Size = 9
$newPart.Paragraphs().ParagraphFormat.Alignment = 2  # ppAlignCenter, consistent with algn="ctr"

# Explicitly resize the shape so its height matches the target cy=630942 EMU
# (note: +49.68051181102362pt, a hair above the exact conversion, avoids
# truncation of the pt->EMU conversion landing 1 EMU short)
$shp.Height = 49.68051181102362

# --- Picture 6 (second occurrence, shape index 24): shift down to y=2696567 EMU ---
$s.Shapes.Item(24).Top = 212.3281496062992

# --- Picture 10 (second occurrence, shape index 25): shift down to y=2700687 EMU ---
$s.Shapes.Item(25).Top = 212.6525590551181

# --- TextBox 23551 (shape index 26): shift down to y=2909697 EMU ---
$s.Shapes.Item(26).Top = 229.11003937007877

# --- TextBox 23554 (shape index 27): shift down to y=2955053 EMU ---
$s.Shapes.Item(27).Top = 232.68137795275592
